$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the anchor paragraph: the blank paragraph that immediately
# follows the "Latency Trend - ..." bullet item (the first of the two
# blank paragraphs before the next page-filler paragraphs).
# ---------------------------------------------------------------------
$anchorRange = $d.Content
$anchorRange.Find.Execute(
    "Latency Trend - Steady until system resource stress kicks in, then increases exponentially.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchorRange.Paragraphs(1)
# The first blank paragraph right after the "Latency Trend" bullet item.
$blank1 = $anchorPara.Next()

$insertPoint = $blank1.Range
$insertPoint.Collapse(0)

# ---------------------------------------------------------------------
# Create 6 new blank paragraphs right after $blank1: one for the
# "Behavioural Summary" heading, three for the bulleted list items, and
# two trailing blank paragraphs.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 6; $i++) {
    $insertPoint.InsertParagraphAfter()
    $insertPoint.Collapse(0)
}

$headingPara = $blank1.Next()
$bullet1Para = $headingPara.Next()
$bullet2Para = $bullet1Para.Next()
$bullet3Para = $bullet2Para.Next()

# ---------------------------------------------------------------------
# "Behavioural Summary" heading - Subtitle style, bold paragraph mark.
# ---------------------------------------------------------------------
$headingPara.Range.Text = "Behavioural Summary"
$headingPara.Style = "Subtitle"
$headingPara.Range.Font.Bold = $true
$headingTextOnly = $d.Range($headingPara.Range.Start, $headingPara.Range.End - 1)
$headingTextOnly.Font.Bold = $false

# ---------------------------------------------------------------------
# Bulleted list - create a brand-new bullet list definition (mirrors
# Word applying "Bullets" to a fresh location) using a throw-away
# paragraph, then transplant just the list (no ListParagraph style) onto
# the three real bullet paragraphs.
# ---------------------------------------------------------------------
$scratchRange = $bullet3Para.Range
$scratchRange.Collapse(0)
$scratchRange.InsertParagraphAfter()
$scratchRange.Collapse(0)
$scratchPara = $bullet3Para.Next()
$scratchPara.Range.Text = "scratch"
$scratchPara.Range.ListFormat.ApplyBulletDefault()
$bulletList = $scratchPara.Range.ListFormat.List

$bullet1Para.Range.Text = "This hash generator service is CPU-heavy, with some memory sensitivity under load."
$bullet1Para.Range.ListFormat.List = $bulletList
$bullet1Para.Range.ListFormat.ListLevelNumber = 1

$bullet2Para.Range.Text = "Latency remained low while CPU and memory were under their thresholds."
$bullet2Para.Range.ListFormat.List = $bulletList
$bullet2Para.Range.ListFormat.ListLevelNumber = 1

$bullet3Para.Range.Text = "Once CPU became saturated and memory started oscillating, latency spiked " + [char]0x2014 + " the system was cascading into performance degradation."
$bullet3Para.Range.ListFormat.List = $bulletList
$bullet3Para.Range.ListFormat.ListLevelNumber = 1

# Remove the scratch paragraph (and its trailing paragraph mark) now
# that the list definition it minted has been reused above.
$scratchPara.Range.Delete()
